$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first two data rows (old rows 2 and 3); remaining rows shift up.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# Append the new accelerometer samples collected on may 9th (rows 20-31).
$ws.Cells.Item(20, 1).Value = 4.854760260603785
$ws.Cells.Item(20, 2).Value = -2.305751679709482
$ws.Cells.Item(20, 3).Value = 28.47659013605825

$ws.Cells.Item(21, 1).Value = 3.906654929683194
$ws.Cells.Item(21, 2).Value = 12.70785129771528
$ws.Cells.Item(21, 3).Value = 6.405859976332387

$ws.Cells.Item(22, 1).Value = 36.5629992182982
$ws.Cells.Item(22, 2).Value = -53.32195605627643
$ws.Cells.Item(22, 3).Value = 10.55735367024117

$ws.Cells.Item(23, 1).Value = 24.49468200994281
$ws.Cells.Item(23, 2).Value = -0.5662195973802113
$ws.Cells.Item(23, 3).Value = -16.59405546058915

$ws.Cells.Item(24, 1).Value = -4.553085171799378
$ws.Cells.Item(24, 2).Value = 1.362173399773924
$ws.Cells.Item(24, 3).Value = -25.30839441912213

$ws.Cells.Item(25, 1).Value = -26.14233835142689
$ws.Cells.Item(25, 2).Value = -26.40326150402269
$ws.Cells.Item(25, 3).Value = 16.80009495834494

$ws.Cells.Item(26, 1).Value = -10.68485793376924
$ws.Cells.Item(26, 2).Value = 1.987670397866637
$ws.Cells.Item(26, 3).Value = 5.468654218302412

$ws.Cells.Item(27, 1).Value = 11.52886270721615
$ws.Cells.Item(27, 2).Value = 5.051807858825063
$ws.Cells.Item(27, 3).Value = 21.48786570044184

$ws.Cells.Item(28, 1).Value = -43.52052723983955
$ws.Cells.Item(28, 2).Value = -14.64106700646965
$ws.Cells.Item(28, 3).Value = -27.88186194993908

$ws.Cells.Item(29, 1).Value = -13.80793043283201
$ws.Cells.Item(29, 2).Value = 3.756517190200164
$ws.Cells.Item(29, 3).Value = -21.9789466857904

$ws.Cells.Item(30, 1).Value = -4.44925512771311
$ws.Cells.Item(30, 2).Value = -2.411408648771987
$ws.Cells.Item(30, 3).Value = -8.979007841774775

$ws.Cells.Item(31, 1).Value = 4.063052345725026
$ws.Cells.Item(31, 2).Value = -3.658945868996846
$ws.Cells.Item(31, 3).Value = 20.40397767459713
